$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("戶號") and give it the new header "撥款".
# This shifts the existing D:J columns (戶號, 額度, 金額, 戶號, 電話, 已繳月份, 撥款日)
# one place to the right, into E:K.
$ws.Columns("D").EntireColumn.Insert()
$ws.Range("D1").Value = "撥款"

# Give the freshly inserted column a sensible (non-default) width like its neighbours.
$ws.Columns("D").ColumnWidth = 5.08

# The worksheet's "Database" defined range grew by the extra column, so extend it
# from column I to column J to keep covering the header row.
$wb.Names.Item("Database").RefersTo = "=LAW7U1Pqp!`$A`$1:`$J`$1"

# Leave the selection on the newly typed header's column, one row down - this is
# where the cursor naturally ends up after typing the new header text.
$ws.Range("D2").Select() | Out-Null
